$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '59.519.22'
$ws.Cells.Item(2, 5).Value = '  -1.87%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.637.05'
$ws.Cells.Item(3, 5).Value = '  -0.29%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.04%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''517.13'
$ws.Cells.Item(5, 5).Value = '  -1.61%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''150.30'
$ws.Cells.Item(6, 5).Value = '  -2.23%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.30%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.06%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.665.86'
$ws.Cells.Item(9, 5).Value = '  +0.31%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''6.56'
$ws.Cells.Item(10, 5).Value = '  +1.44%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -1.24%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -1.86%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -0.83%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '3.104.04'
$ws.Cells.Item(14, 5).Value = '  -0.19%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '59.263.35'
$ws.Cells.Item(15, 5).Value = '  -2.34%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -1.75%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -0.80%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.663.02'
$ws.Cells.Item(18, 5).Value = '  +0.50%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -2.09%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''348.98'
$ws.Cells.Item(20, 5).Value = '  -0.82%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''10.65'
$ws.Cells.Item(21, 5).Value = '  +0.54%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''6.25'
$ws.Cells.Item(22, 5).Value = '  -0.10%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''0.999'
$ws.Cells.Item(23, 5).Value = '  -0.13%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''61.46'
$ws.Cells.Item(24, 5).Value = '  +0.37%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''0.428'
$ws.Cells.Item(25, 5).Value = '  +0.28%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '2.754.46'
$ws.Cells.Item(26, 5).Value = '  -0.47%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Kaspa'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(27, 4).Value = '''0.162'
$ws.Cells.Item(27, 5).Value = '  -2.58%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(28, 4).Value = '''0.991'
$ws.Cells.Item(28, 5).Value = '  -0.70%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -1.14%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''7.18'
$ws.Cells.Item(30, 5).Value = '  -0.80%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -0.22%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''6.44'
$ws.Cells.Item(32, 5).Value = '  +4.81%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''19.21'
$ws.Cells.Item(33, 5).Value = '  -0.60%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -2.43%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''149.56'
$ws.Cells.Item(35, 5).Value = '  -0.14%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +17.50%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +0.65%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''1.17'
$ws.Cells.Item(38, 5).Value = '  -1.38%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''0.895'
$ws.Cells.Item(39, 5).Value = '  -0.84%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''36.48'
$ws.Cells.Item(40, 5).Value = '  -0.86%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''1.46'
$ws.Cells.Item(41, 5).Value = '  -0.32%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.54%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''293.60'
$ws.Cells.Item(43, 5).Value = '  -4.02%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''0.629'
$ws.Cells.Item(44, 5).Value = '  -0.99%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.93%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''0.993'
$ws.Cells.Item(46, 5).Value = '  -0.56%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''19.78'
$ws.Cells.Item(47, 5).Value = '  -1.50%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -1.59%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''4.84'
$ws.Cells.Item(49, 5).Value = '  +0.11%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -2.06%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''19.07'
$ws.Cells.Item(51, 5).Value = '  +0.17%  '
